$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: convert SmartScore inline-string cells to numeric ---
$ws.Range("G2").Value = 0.627
$ws.Range("J2").Value = 0.469
$ws.Range("M2").Value = 0.458
$ws.Range("P2").Value = 0.696
$ws.Range("S2").Value = 0.64
$ws.Range("V2").Value = 0.577
$ws.Range("Y2").Value = 0.679
$ws.Range("AB2").Value = 0.545
$ws.Range("AE2").Value = 0.516

# --- Row 3: new record for Juan Luis ---
$ws.Range("A3").Value = 'Juan Luis'
$ws.Range("B3").Value = 24
$ws.Range("C3").Value = 'Masculino'
$ws.Range("D3").Value = '2025-10-28 05:27:34'
$ws.Range("E3").Value = '{
  "portion": 0.8,
  "diet": 1.0,
  "salt": 0.8,
  "fat": 0.0,
  "natural": 0.8,
  "convenience": 0.2,
  "price": 0.2
}'
$ws.Range("F3").Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range("G3").Value = '''0.572'
$ws.Range("H3").Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("I3").Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range("J3").Value = '''0.514'
$ws.Range("K3").Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("L3").Value = 'Nongshim Shin Ramyun'
$ws.Range("M3").Value = '''0.409'
$ws.Range("N3").Value = 'Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio'
$ws.Range("O3").Value = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range("P3").Value = '''0.845'
$ws.Range("Q3").Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range("R3").Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range("S3").Value = '''0.618'
$ws.Range("T3").Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("U3").Value = 'Annie’s Shells & White Cheddar'
$ws.Range("V3").Value = '''0.602'
$ws.Range("W3").Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range("X3").Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range("Y3").Value = '''0.769'
$ws.Range("Z3").Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range("AA3").Value = 'Kitchens of India Variety Pack'
$ws.Range("AB3").Value = '''0.503'
$ws.Range("AC3").Value = 'Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad'
$ws.Range("AD3").Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range("AE3").Value = '''0.423'
$ws.Range("AF3").Value = 'Portátil, saludable, fácil, buena textura, sabor suave'
